# The deck shipped with two themes: theme1.xml ("Integral", wired to the
# slide master -> all slides) and theme2.xml ("Office Theme", wired to the
# notes master). The commit swaps the two themes so the slide master (and
# therefore every slide) now carries the "Office Theme" palette while the
# notes master keeps the "Integral" one.
#
# The PowerPoint object model exposes theme colors per slide/master through
# ThemeColorScheme (12 slots: dk1, lt1, dk2, lt2, accent1-6, hlink,
# folHlink) - font scheme and format scheme (fills/lines/effects) are
# identical between the two themes here, so re-pointing the 12 colors is
# enough to reproduce the swap.

$p = $ppt.ActivePresentation

# Office Theme palette (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink)
# expressed as COM RGB long values (R + G*256 + B*65536).
$officeThemeColors = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

# All slides share the one slide master / theme (theme1.xml), so applying
# the new palette through the first slide's ThemeColorScheme updates the
# shared theme part for the whole deck.
$s = $p.Slides.Item(1)
for ($i = 1; $i -le $officeThemeColors.Count; $i++) {
    $s.ThemeColorScheme.Item($i).RGB = $officeThemeColors[$i - 1]
}
